$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.452.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5204"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4359"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.29"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +16.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08829"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.151"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.26"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.098.29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.664"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.651"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001119"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06589"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.23"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.253"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.506.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.335.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.24"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.550"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.28"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.182"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1065"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.637"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.132"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.903"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.16"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02571"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06814"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.65"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.427"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2247"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6858"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.261"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.45%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.37%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6340"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.188"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.621"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.13%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.239"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.04%  "
